# Scene.xlsx configuration update:
#  - Row 2 (ID=0, scene "Login") is repurposed into a new "clone" scene entry
#    that points at the CloneScene folder and gets ID 3.
#  - Rows 3 and 4 (previously SceneName "Stage001") are renamed to "newscene".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FilePath, ID and SceneName change (ID set first so shared-string
# ordering mirrors the authored workbook).
$ws.Range("B2").Value = "3"
$ws.Range("A2").Value = "../../NFDataCfg/Ini/NFZoneServer/Scene/CloneScene/"
$ws.Range("F2").Value = "clone"

# Row 3: SceneName changes from Stage001 to newscene.
$ws.Range("F3").Value = "newscene"

# Row 4: SceneName changes from Stage001 to newscene.
$ws.Range("F4").Value = "newscene"

# Update the active selection left behind in the sheet view.
[void]$ws.Range("H8").Select()
